# CMMC Fledge Milestone 1 Report - Final edit pass
# Applies the textual revisions + two GitHub hyperlinks described by the
# commit "Milestone 1 Report Final".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $find"
    }
    return $ok
}

# 1) "...process of Cybersecurity Maturity Model Certification (CMMC) compliance."
#    -> "...process of the Cybersecurity Maturity Model Certification (CMMC)."
Replace-Text "process of Cybersecurity Maturity Model Certification (CMMC) compliance." "process of the Cybersecurity Maturity Model Certification (CMMC)."

# 2) Tighten the "figuring out the best way to support..." sentence.
Replace-Text "research into CMMC documentation and figuring out the best way to support the customers of the CMMC Fledge System. By achieving a greater understanding of the process and some of its nuances," "research into CMMC documentation. By achieving this greater understanding of the process and some of its nuances,"

# 3) Drop "exactly" from "...CMMC Fledge System exactly is can be seen."
Replace-Text "what the CMMC Fledge System exactly is can be seen." "what the CMMC Fledge System is can be seen."

# 4) Drop "itself" from "...CMMC Fledge System itself."  (before the GitHub hyperlink is created)
Replace-Text "A GitHub repository has been established to hold the documents drafted thus far and the future implementation (code, files, etc.) of the CMMC Fledge System itself." "A GitHub repository has been established to hold the documents drafted thus far and the future implementation (code, files, etc.) of the CMMC Fledge System."

# 5) Remove "This system mainly consists of an interactive webpage and a database. "
Replace-Text "This milestone set the foundation for the development of the CMMC Fledge System. This system mainly consists of an interactive webpage and a database. Between" "This milestone set the foundation for the development of the CMMC Fledge System. Between"

# 6) Turn the "GitHub repository" mention into a hyperlink.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("GitHub repository")
if ($found1) {
    $d.Hyperlinks.Add($rng1, "https://github.com/CMMCFledge/CMMCFledge") | Out-Null
} else {
    Write-Output "WARNING: GitHub repository text not found for hyperlink"
}

# 7) Turn the "GitHub" mention (before " commits") into a hyperlink.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("GitHub commits")
if ($found2) {
    $ghEnd = $rng2.Start + 6
    $rngGh = $d.Range($rng2.Start, $ghEnd)
    $d.Hyperlinks.Add($rngGh, "https://github.com/CMMCFledge") | Out-Null
} else {
    Write-Output "WARNING: GitHub commits text not found for hyperlink"
}

Write-Output "Edits applied"
